$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")
$ws.Range("A7").Value = "Notifications"
$ws.Range("A7").Select()
